$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must remain TEXT even when it looks like a
# plain number (e.g. "524.79"). Excel auto-converts Range.Value = "524.79"
# to a numeric cell, which would also diverge from the original inline-string
# typing. Routing the literal through a scratch formula cell and pasting
# values-only keeps the destination cell's existing (default) style intact
# while forcing genuine text content.
function Set-TextValue($addr, $val) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $val + '"'
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "58.541.54"
$ws.Range("E2").Value = "  +0.50%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.539.91"
$ws.Range("E3").Value = "  +2.14%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.10%  "

# Row 5 - BNB
Set-TextValue "D5" "524.79"
$ws.Range("E5").Value = "  +0.71%  "

# Row 6 - Solana
Set-TextValue "D6" "133.56"
$ws.Range("E6").Value = "  -0.83%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.41%  "

# Row 8 - XRP
Set-TextValue "D8" "0.565"
$ws.Range("E8").Value = "  +1.02%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.537.34"
$ws.Range("E9").Value = "  +1.33%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.64%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -1.40%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  -2.46%  "

# Row 13 - Cardano
Set-TextValue "D13" "0.333"
$ws.Range("E13").Value = "  -1.99%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.988.04"
$ws.Range("E14").Value = "  +2.04%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "58.506.46"
$ws.Range("E15").Value = "  +0.58%  "

# Row 16 - Avalanche
Set-TextValue "D16" "22.26"
$ws.Range("E16").Value = "  +0.51%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -0.18%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.537.48"
$ws.Range("E18").Value = "  +1.68%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  +0.12%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "323.65"
$ws.Range("E20").Value = "  +0.57%  "

# Row 21 - Polkadot
$ws.Range("E21").Value = "  -0.15%  "

# Row 22 - Uniswap
Set-TextValue "D22" "6.17"
$ws.Range("E22").Value = "  +6.89%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.12%  "

# Row 24 - Litecoin
Set-TextValue "D24" "65.17"
$ws.Range("E24").Value = "  +0.64%  "

# Row 25 - Polygon
$ws.Range("E25").Value = "  -1.07%  "

# Row 26 - Binance-PegBSC-USD
Set-TextValue "D26" "0.998"
$ws.Range("E26").Value = "  +0.61%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  -1.06%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("E28").Value = "  -0.22%  "

# Row 29 - PEPE
$ws.Range("E29").Value = "  +0.19%  "

# Row 30 - PancakeSwap
Set-TextValue "D30" "1.73"
$ws.Range("E30").Value = "  +1.74%  "

# Row 31 - Monero
Set-TextValue "D31" "168.46"
$ws.Range("E31").Value = "  -0.63%  "

# Row 32 - Fetch.AI
Set-TextValue "D32" "1.20"
$ws.Range("E32").Value = "  +0.98%  "

# Row 33 - Aptos
Set-TextValue "D33" "6.32"
$ws.Range("E33").Value = "  -0.50%  "

# Row 35 - FirstDigitalUSD
Set-TextValue "D35" "0.997"
$ws.Range("E35").Value = "  +0.13%  "

# Row 36 - EthereumClassic
$ws.Range("E36").Value = "  +0.92%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  -4.66%  "

# Row 38 - NEARProtocol
$ws.Range("E38").Value = "  -2.26%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +0.96%  "

# Row 40 - OKB
Set-TextValue "D40" "36.48"
$ws.Range("E40").Value = "  -0.51%  "

# Row 41 - SuiNetwork
Set-TextValue "D41" "0.777"
$ws.Range("E41").Value = "  -2.93%  "

# Row 42 - Bittensor
Set-TextValue "D42" "278.86"
$ws.Range("E42").Value = "  +0.93%  "

# Row 43 - Filecoin
$ws.Range("E43").Value = "  +0.23%  "

# Row 44 - RenderToken
Set-TextValue "D44" "5.04"
$ws.Range("E44").Value = "  -0.63%  "

# Row 45 - was Mantle, becomes Aave (swapped with row 46)
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D45" "130.20"
$ws.Range("E45").Value = "  +4.84%  "

# Row 46 - was Aave, becomes Mantle (swapped with row 45)
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D46" "0.604"
$ws.Range("E46").Value = "  +0.84%  "

# Row 47 - Stellar
$ws.Range("E47").Value = "  +0.81%  "

# Row 48 - Hedera
Set-TextValue "D48" "0.0502"
$ws.Range("E48").Value = "  +1.90%  "

# Row 49 - EnergySwap
$ws.Range("E49").Value = "  +0.25%  "

# Row 50 - VeChain
$ws.Range("E50").Value = "  +0.12%  "

# Row 51 - InjectiveProtocol
Set-TextValue "D51" "17.07"
$ws.Range("E51").Value = "  -0.62%  "

# Clean up the scratch cell so no stray content / used-range growth remains.
$ws.Range("ZZ1").Clear()
